# Commit: Created UI for session setup panel config for Time mode interval.
# Adds 15 new rows (39-53) to the "Translation" sheet describing new
# translation text entries for the Time-mode session setup panel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Each entry: row number, TEXT ID (col B), TYPOGRAPHY NAME (col C),
# ALIGNMENT (col D), DIRECTION (col E), GB text (col F)
$rows = @(
    @{ Row = 39; B = "SingleUseId36"; C = "Default"; D = "Left";   E = "LTR"; F = "TI MAX RANGE" },
    @{ Row = 40; B = "SingleUseId37"; C = "Default"; D = "Center"; E = "LTR"; F = "Single" },
    @{ Row = 41; B = "SingleUseId39"; C = "Default"; D = "Left";   E = "LTR"; F = "Meas Rate" },
    @{ Row = 42; B = "SingleUseId40"; C = "Default"; D = "Center"; E = "LTR"; F = "<value> ms" },
    @{ Row = 43; B = "SingleUseId41"; C = "Default"; D = "Center"; E = "LTR"; F = "Continuous" },
    @{ Row = 44; B = "SingleUseId42"; C = "Default"; D = "Center"; E = "LTR"; F = "<value>" },
    @{ Row = 45; B = "SingleUseId43"; C = "Default"; D = "Left";   E = "LTR"; F = "0" },
    @{ Row = 46; B = "SingleUseId44"; C = "Default"; D = "Left";   E = "LTR"; F = "0" },
    @{ Row = 47; B = "SingleUseId45"; C = "Default"; D = "Left";   E = "LTR"; F = "Stamps Number" },
    @{ Row = 48; B = "SingleUseId46"; C = "Large";   D = "Left";   E = "LTR"; F = "X" },
    @{ Row = 49; B = "SingleUseId47"; C = "Default"; D = "Left";   E = "LTR"; F = "Repeat" },
    @{ Row = 50; B = "SingleUseId48"; C = "Default"; D = "Center"; E = "LTR"; F = "<value>" },
    @{ Row = 51; B = "SingleUseId49"; C = "Default"; D = "Left";   E = "LTR"; F = "0" },
    @{ Row = 52; B = "SingleUseId50"; C = "Default"; D = "Center"; E = "LTR"; F = "<value> s" },
    @{ Row = 53; B = "SingleUseId51"; C = "Default"; D = "Left";   E = "LTR"; F = "0" }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E

    $fCell = $ws.Cells.Item($r, 6)
    $fText = $entry.F

    $looksNumeric = $fText -match '^-?\d+(\.\d+)?$'

    if ($looksNumeric) {
        # Writing a bare numeric-looking literal (e.g. "0") via .Value would be
        # auto-coerced to a Number by Excel's normal type inference, and using
        # a leading apostrophe to force text adds a quotePrefix cell style that
        # the source workbook does not have. Going through a text formula and
        # collapsing it to a value via Copy/PasteSpecial(xlPasteValues) keeps
        # the result a plain shared-string text cell with no extra styling.
        $fCell.Formula = '="' + $fText + '"'
        $fCell.Copy()
        $fCell.PasteSpecial(-4163) # xlPasteValues
    }
    else {
        $fCell.Value = $fText
    }
}

$excel.CutCopyMode = 0
